# Apply weekly fruit/vegetable price update: swap the data values between
# row 2 and row 3 for columns D (Fecha), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado) and S (Precio $/Kg).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D - Fecha (numeric date serials, keep raw numbers so format stays intact)
$ws.Range("D2").Value = 44217
$ws.Range("D3").Value = 44209

# Column M - Volumen
$ws.Range("M2").Value = 200
$ws.Range("M3").Value = 100

# Column N - Precio minimo
$ws.Range("N2").Value = 11000
$ws.Range("N3").Value = 10000

# Column O - Precio maximo
$ws.Range("O2").Value = 12000
$ws.Range("O3").Value = 11000

# Column P - Precio promedio ponderado
$ws.Range("P2").Value = 11500
$ws.Range("P3").Value = 10500

# Column S - Precio $/Kg
$ws.Range("S2").Value = 821
$ws.Range("S3").Value = 750
